# Auto update Excel log
# Appends new sensor-log rows to the three mmWave sheets, matching the
# upstream "Auto update Excel log" commit.
#
# NOTE: the "Date" column (A) holds plain text like "2026-02-01". If we
# assign that directly via .Value, Excel's COM layer auto-detects it as a
# date and stores a date serial number instead of text. To preserve the
# original text semantics we briefly force the cell to a text format
# before writing the value, then restore the cell style to "Normal" so no
# leftover number-format styling is left on the cell (matches the
# plain/unstyled cells used throughout the log).

function Set-LogRow {
    param($ws, $r, $date, $timestamp, $hour, $location, $value, $status)

    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $date
    $ws.Cells.Item($r, 1).Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $timestamp
    $ws.Cells.Item($r, 3).Value = $hour
    $ws.Cells.Item($r, 4).Value = $location
    $ws.Cells.Item($r, 5).Value = $value
    $ws.Cells.Item($r, 6).Value = $status
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# mmWave(InBed): rows 123-130 (E column is a text status: "In Bed")
# ---------------------------------------------------------------------
$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")

Set-LogRow $wsInBed 123 "2026-02-01" "20:37:33" "20:00" "Bedroom" "In Bed" "Occupied"
Set-LogRow $wsInBed 124 "2026-02-01" "20:37:34" "20:00" "Bedroom" "In Bed" "Occupied"
Set-LogRow $wsInBed 125 "2026-02-01" "20:37:36" "20:00" "Bedroom" "In Bed" "Occupied"
Set-LogRow $wsInBed 126 "2026-02-01" "20:37:38" "20:00" "Bedroom" "In Bed" "Occupied"
Set-LogRow $wsInBed 127 "2026-02-01" "20:37:40" "20:00" "Bedroom" "In Bed" "Occupied"
Set-LogRow $wsInBed 128 "2026-02-01" "20:37:42" "20:00" "Bedroom" "In Bed" "Occupied"
Set-LogRow $wsInBed 129 "2026-02-01" "20:37:44" "20:00" "Bedroom" "In Bed" "Occupied"
Set-LogRow $wsInBed 130 "2026-02-01" "20:37:46" "20:00" "Bedroom" "In Bed" "Occupied"

# ---------------------------------------------------------------------
# mmWave(BR): rows 117-123 (E column is a numeric value)
# ---------------------------------------------------------------------
$wsBR = $wb.Worksheets.Item("mmWave(BR)")

Set-LogRow $wsBR 117 "2026-02-01" "20:37:36" "20:00" "Bedroom" 29 "Occupied"
Set-LogRow $wsBR 118 "2026-02-01" "20:37:38" "20:00" "Bedroom" 11 "Occupied"
Set-LogRow $wsBR 119 "2026-02-01" "20:37:40" "20:00" "Bedroom" 34 "Occupied"
Set-LogRow $wsBR 120 "2026-02-01" "20:37:41" "20:00" "Bedroom" 36 "Occupied"
Set-LogRow $wsBR 121 "2026-02-01" "20:37:44" "20:00" "Bedroom" 60 "Occupied"
Set-LogRow $wsBR 122 "2026-02-01" "20:37:46" "20:00" "Bedroom" 2 "Occupied"
Set-LogRow $wsBR 123 "2026-02-01" "20:37:48" "20:00" "Bedroom" 1 "Occupied"

# ---------------------------------------------------------------------
# mmWave(HR): rows 117-123 (E column is a numeric value)
# ---------------------------------------------------------------------
$wsHR = $wb.Worksheets.Item("mmWave(HR)")

Set-LogRow $wsHR 117 "2026-02-01" "20:37:35" "20:00" "Bedroom" 77 "Occupied"
Set-LogRow $wsHR 118 "2026-02-01" "20:37:37" "20:00" "Bedroom" 59 "Occupied"
Set-LogRow $wsHR 119 "2026-02-01" "20:37:39" "20:00" "Bedroom" 82 "Occupied"
Set-LogRow $wsHR 120 "2026-02-01" "20:37:41" "20:00" "Bedroom" 84 "Occupied"
Set-LogRow $wsHR 121 "2026-02-01" "20:37:43" "20:00" "Bedroom" 108 "Occupied"
Set-LogRow $wsHR 122 "2026-02-01" "20:37:45" "20:00" "Bedroom" 50 "Occupied"
Set-LogRow $wsHR 123 "2026-02-01" "20:37:47" "20:00" "Bedroom" 49 "Occupied"
